# feat: add 2022-Q4 data
#
# The workbook currently has two sheets: "总计" (totals) and "2022-Q3"
# (fund holdings detail for 2022-Q3). This script:
#   1. Duplicates the "2022-Q3" sheet (so its original data survives
#      unchanged) right after itself -> becomes the new "2022-Q3" sheet.
#   2. Renames the original "2022-Q3" sheet to "2022-Q4" and overwrites
#      its data with the new Q4 fund-holdings detail.
#   3. Inserts a new summary row for 2022-Q4 at the top of the "总计"
#      sheet (row 2), pushing the existing 2022-Q3 summary row down to
#      row 3 and bumping its running index from 0 to 1.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

# --- Step 1: duplicate the 2022-Q3 sheet, placing the copy right after it.
# The copy keeps all of the original 2022-Q3 data/format untouched and
# becomes the new, permanent "2022-Q3" sheet.
$q3Sheet.Copy($null, $q3Sheet)
$q3Copy = $wb.Worksheets.Item(3)

# --- Step 2: turn the original sheet into the 2022-Q4 sheet (rename it
# out of the way first so the copy can reclaim the "2022-Q3" name).
$q4Sheet = $q3Sheet
$q4Sheet.Name = "2022-Q4"
$q3Copy.Name = "2022-Q3"

$q4Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q4Headers.Length; $i++) {
    $cell = $q4Sheet.Cells.Item(1, $i + 2)
    $cell.Value = $q4Headers[$i]
    $cell.Style = $totalSheet.Cells.Item(1, 2).Style
}

$q4Data = @(
    @("011201", "财通优势行业轮动混合A", "8.38", "88.19", "4.21", "0.3528", 7),
    @("501085", "财通科创主题灵活配置混合（LOF）", "4.19", "89.51", "6.56", "0.2749", 3),
    @("010874", "泰康品质生活混合A", "6.77", "84.20", "3.18", "0.2153", 7),
    @("010994", "博时创新经济混合A", "3.80", "86.65", "5.46", "0.2075", 2),
    @("013680", "华安品质甄选混合A", "12.95", "73.22", "1.23", "0.1593", 10),
    @("010875", "泰康品质生活混合C", "3.33", "84.20", "3.18", "0.1059", 7),
    @("010995", "博时创新经济混合C", "1.71", "86.65", "5.46", "0.0934", 2),
    @("016336", "博时卓远成长一年持有期股票A", "1.63", "51.15", "5.31", "0.0866", 2),
    @("013681", "华安品质甄选混合C", "5.10", "73.22", "1.23", "0.0627", 10),
    @("501001", "财通多策略精选混合（LOF）", "0.74", "87.17", "4.32", "0.0320", 3),
    @("016337", "博时卓远成长一年持有期股票C", "0.46", "51.15", "5.31", "0.0244", 2),
    @("011202", "财通优势行业轮动混合C", "0.34", "88.19", "4.21", "0.0143", 7),
    @("003938", "南方荣尊混合A", "0.21", "29.78", "1.89", "0.0040", 6),
    @("003939", "南方荣尊混合C", "0.09", "29.78", "1.89", "0.0017", 6)
)

for ($r = 0; $r -lt $q4Data.Length; $r++) {
    $row = $r + 2
    $rec = $q4Data[$r]

    $idxCell = $q4Sheet.Cells.Item($row, 1)
    $idxCell.Value = $r
    $idxCell.Style = $totalSheet.Cells.Item(2, 1).Style

    # Columns B:G are textual (fund code keeps leading zeros, the
    # figures keep their original decimal-string formatting) rather
    # than being auto-coerced to numbers by Excel.
    for ($col = 2; $col -le 7; $col++) {
        $dataCell = $q4Sheet.Cells.Item($row, $col)
        $dataCell.NumberFormat = "@"
        $dataCell.Value = $rec[$col - 2]
    }

    $q4Sheet.Cells.Item($row, 8).Value = $rec[6]
}

# --- Step 3: insert the 2022-Q4 summary row above the existing 2022-Q3
# summary row on the "总计" sheet, then update the 2022-Q3 row's index.
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 1).Style = $totalSheet.Cells.Item(3, 1).Style
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 14
$totalSheet.Cells.Item(2, 4).Value = 1.63

$totalSheet.Cells.Item(3, 1).Value = 1
